# Add two new "40p jumper wire" line items (rows 26 and 27) to the
# "부품 리스트" (parts list) sheet, matching the row-25 formatting, and
# update the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "부품 리스트"

# --- Copy the formatting of the last existing data row (25) down onto
#     the two new rows (26, 27) so columns A-E pick up the same styles
#     (centered number style for A/B/C, the summed style for D, and the
#     Hyperlink style for E). ---
$ws.Range("A25:E25").Copy() | Out-Null
$ws.Range("A26:E27").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 26: "40p 암암 점퍼선" ---
$url26 = "http://www.11st.co.kr/product/SellerProductDetail.tmall?method=getSellerProductDetail&prdNo=1323271899&trTypeCd=PW53&trCtgrNo=585021"
$ws.Hyperlinks.Add($ws.Range("E26"), $url26) | Out-Null
$ws.Range("A26").Value = "40p 암암 점퍼선"
$ws.Range("B26").Value = 1100
$ws.Range("C26").Value = 1
# D26 keeps the shared formula (=B26*C26) inherited from the paste above.

# --- Row 27: "40p 수수 점퍼선" ---
$url27 = "http://www.11st.co.kr/product/SellerProductDetail.tmall?method=getSellerProductDetail&prdNo=2579440230&trTypeCd=PW24&trCtgrNo=585021&lCtgrNo=1001362&mCtgrNo=1002187"
$ws.Hyperlinks.Add($ws.Range("E27"), $url27) | Out-Null
$ws.Range("A27").Value = "40p 수수 점퍼선"
$ws.Range("B27").Value = 1700
$ws.Range("C27").Value = 1
# D27 keeps the shared formula (=B27*C27) inherited from the paste above.

# Re-apply the row-25 formats once more: adding the hyperlinks above can
# nudge the cell style of E26/E27 onto a freshly-created (but equivalent)
# style record, so re-copying the formats collapses them back onto the
# existing "하이퍼링크" style used throughout column E.
$ws.Range("A25:E25").Copy() | Out-Null
$ws.Range("A26:E27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# D30 = SUM(D2:D29) recalculates automatically to include the two new rows.

# --- Update the sheet view: selection moves to C31 and the view scrolls
#     down so row 16 is at the top. ---
$ws.Activate()
$ws.Range("C31").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
